# Update weekly excess mortality analyses
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct a handful of already-reported weekly death counts (column G)
#     Column I (Oversterfte = G - H) recalculates automatically since it holds
#     a formula in every row.
$ws.Range("G8").Value  = 4305
$ws.Range("G21").Value = 2528
$ws.Range("G23").Value = 2667
$ws.Range("G24").Value = 2640
$ws.Range("G26").Value = 2853
$ws.Range("G31").Value = 2891
$ws.Range("G33").Value = 3019
$ws.Range("G34").Value = 3212
$ws.Range("G35").Value = 3444
$ws.Range("G36").Value = 3674
$ws.Range("G37").Value = 3587
$ws.Range("G38").Value = 3552
$ws.Range("G39").Value = 3315
$ws.Range("G40").Value = 3373

# --- 2. Insert a blank separator row just above the totals row, pushing the
#     totals (currently row 42) down to row 43.
$ws.Rows("42:42").Insert()

# --- 3. Fill in the newly reported week 49 figures on row 41.
$ws.Range("F41").Value = 49
$ws.Range("G41").Value = 3448
$ws.Range("H41").Value = 3037
$ws.Range("I41").Formula = "=G41-H41"

# --- 4. Keep the selection where the author left it.
$ws.Range("J41").Select()
